$wb = $excel.ActiveWorkbook

# --- Typography sheet: move "Widget Wildcard Characters" from the
#     "Clock" typography (row 9) to the "Default" typography (row 4) ---
$typography = $wb.Worksheets.Item("Typography")
$typography.Range("H4").Value = "0123456789 :APM"
$typography.Range("H9").Value = ""

# --- Translation sheet: remove the old "Clock" row, shift subsequent
#     rows up, and append a new "Digital_clock" row at the bottom ---
$translation = $wb.Worksheets.Item("Translation")

$translation.Range("B11").Value = "SingleUseId9"
$translation.Range("C11").Value = "Default"
$translation.Range("D11").Value = "Center"
$translation.Range("E11").Value = "LTR"
$translation.Range("F11").Value = "Cancel"

$translation.Range("B12").Value = "SingleUseId10"
$translation.Range("C12").Value = "Default"
$translation.Range("D12").Value = "Left"
$translation.Range("E12").Value = "LTR"
$translation.Range("F12").Value = "Ok"

$translation.Range("B13").Value = "SingleUseId11"
$translation.Range("C13").Value = "modalWindowTitle"
$translation.Range("D13").Value = "Left"
$translation.Range("E13").Value = "LTR"
$translation.Range("F13").Value = "<value>"

$translation.Range("B14").Value = "SingleUseId12"
$translation.Range("C14").Value = "Default"
$translation.Range("D14").Value = "Left"
$translation.Range("E14").Value = "LTR"
$translation.Range("F14").Value = "Password"

$translation.Range("B15").Value = "SingleUseId13"
$translation.Range("C15").Value = "Large"
$translation.Range("D15").Value = "Left"
$translation.Range("E15").Value = "LTR"
$translation.Range("F15").Value = "<value>"

$translation.Range("B16").Value = "Digital_clock"
$translation.Range("C16").Value = "Default"
$translation.Range("D16").Value = "Left"
$translation.Range("E16").Value = "LTR"
$translation.Range("F16").Value = "<time>"
